$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 12 (rows 12-50 shift down to 14-52)
$ws.Rows(12).Insert()
$ws.Rows(12).Insert()

# New row 12 data
$ws.Cells.Item(12,1).Value = 2
$ws.Cells.Item(12,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(12,3).Value = "Coquimbo"
$ws.Cells.Item(12,4).Value = 45274
$ws.Cells.Item(12,5).Value = 4
$ws.Cells.Item(12,6).Value = "Fruta"
$ws.Cells.Item(12,7).Value = 100103
$ws.Cells.Item(12,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(12,9).Value = 100103003
$ws.Cells.Item(12,10).Value = "Damasco"
$ws.Cells.Item(12,11).Value = "Castle Brite"
$ws.Cells.Item(12,12).Value = "Especial"
$ws.Cells.Item(12,13).Value = 100
$ws.Cells.Item(12,14).Value = 26000
$ws.Cells.Item(12,15).Value = 27000
$ws.Cells.Item(12,16).Value = 26500
$ws.Cells.Item(12,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(12,18).Value = "Región de O'Higgins"
$ws.Cells.Item(12,19).Value = 1472
$ws.Cells.Item(12,20).Value = 18

# New row 13 data
$ws.Cells.Item(13,1).Value = 2
$ws.Cells.Item(13,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13,3).Value = "Coquimbo"
$ws.Cells.Item(13,4).Value = 45274
$ws.Cells.Item(13,5).Value = 4
$ws.Cells.Item(13,6).Value = "Fruta"
$ws.Cells.Item(13,7).Value = 100103
$ws.Cells.Item(13,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13,9).Value = 100103003
$ws.Cells.Item(13,10).Value = "Damasco"
$ws.Cells.Item(13,11).Value = "Castle Brite"
$ws.Cells.Item(13,12).Value = "Primera"
$ws.Cells.Item(13,13).Value = 100
$ws.Cells.Item(13,14).Value = 24000
$ws.Cells.Item(13,15).Value = 25000
$ws.Cells.Item(13,16).Value = 24500
$ws.Cells.Item(13,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(13,18).Value = "Región de O'Higgins"
$ws.Cells.Item(13,19).Value = 1361
$ws.Cells.Item(13,20).Value = 18
